$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "sitenumber"

$ws.Range("A2").Value = 36.015714
$ws.Range("B2").Value = -90.799477
$ws.Range("A3").Value = 40.704774
$ws.Range("B3").Value = -73.94186000000001
$ws.Range("A4").Value = 41.505144
$ws.Range("B4").Value = -94.518159
$ws.Range("A5").Value = 46.916976
$ws.Range("B5").Value = -97.996522
$ws.Range("A6").Value = 38.41666
$ws.Range("B6").Value = -122.8662
$ws.Range("A7").Value = 31.127656
$ws.Range("B7").Value = -97.86003599999999
$ws.Range("A8").Value = 33.722
$ws.Range("B8").Value = -117.91413
$ws.Range("A9").Value = 34.04695
$ws.Range("B9").Value = -118.03132
$ws.Range("A10").Value = 29.904808
$ws.Range("B10").Value = -90.136928
$ws.Range("A11").Value = 40.10397
$ws.Range("B11").Value = -83.77349
$ws.Range("A12").Value = 32.9296
$ws.Range("B12").Value = -96.323607
$ws.Range("A13").Value = 40.28939
$ws.Range("B13").Value = -73.98488999999999
$ws.Range("A14").Value = 39.08575
$ws.Range("B14").Value = -89.39008
$ws.Range("A15").Value = 30.314791
$ws.Range("B15").Value = -95.38846700000001
$ws.Range("A16").Value = 35.101556
$ws.Range("B16").Value = -117.991517
$ws.Range("A17").Value = 42.18712
$ws.Range("B17").Value = -83.76461999999999
$ws.Range("A18").Value = 40.767653
$ws.Range("B18").Value = -99.702
$ws.Range("A19").Value = 30.680797
$ws.Range("B19").Value = -96.318895
$ws.Range("A20").Value = 41.71859
$ws.Range("B20").Value = -87.81511999999999
$ws.Range("A21").Value = 39.95076
$ws.Range("B21").Value = -84.99476
$ws.Range("A22").Value = 39.38928
$ws.Range("B22").Value = -74.54725999999999
$ws.Range("A23").Value = 35.08738
$ws.Range("B23").Value = -106.5568
$ws.Range("A24").Value = 34.488866
$ws.Range("B24").Value = -88.207067
$ws.Range("A25").Value = 32.62801
$ws.Range("B25").Value = -96.90831
$ws.Range("A26").Value = 41.67021
$ws.Range("B26").Value = -97.04919
$ws.Range("A27").Value = 36.91443
$ws.Range("B27").Value = -98.9902
$ws.Range("A28").Value = 42.87932
$ws.Range("B28").Value = -71.57973
$ws.Range("A29").Value = 43.05528
$ws.Range("B29").Value = -96.28747
$ws.Range("A30").Value = 39.31225
$ws.Range("B30").Value = -74.59278
$ws.Range("A31").Value = 38.0854
$ws.Range("B31").Value = -122.2573
$ws.Range("A32").Value = 39.95654
$ws.Range("B32").Value = -75.25861999999999
$ws.Range("A33").Value = 40.5612
$ws.Range("B33").Value = -75.5235
$ws.Range("A34").Value = 42.00607
$ws.Range("B34").Value = -76.55620999999999
$ws.Range("A35").Value = 38.573208
$ws.Range("B35").Value = -121.257761
$ws.Range("A36").Value = 39.326667
$ws.Range("B36").Value = -120.182222
$ws.Range("A37").Value = 40.820307
$ws.Range("B37").Value = -73.899063
$ws.Range("A38").Value = 45.01415
$ws.Range("B38").Value = -93.17711
$ws.Range("A39").Value = 33.693743
$ws.Range("B39").Value = -117.954269
$ws.Range("A40").Value = 38.64312
$ws.Range("B40").Value = -90.31562
$ws.Range("A41").Value = 38.258
$ws.Range("B41").Value = -77.1447
$ws.Range("A42").Value = 41.23908
$ws.Range("B42").Value = -85.85663
$ws.Range("A43").Value = 29.563333
$ws.Range("B43").Value = -104.33944
$ws.Range("A44").Value = 37.30398
$ws.Range("B44").Value = -120.47281
$ws.Range("A45").Value = 38.88487
$ws.Range("B45").Value = -77.90571
$ws.Range("A46").Value = 47.7512
$ws.Range("B46").Value = -121.465
$ws.Range("A47").Value = 40.725253
$ws.Range("B47").Value = -74.18385600000001
$ws.Range("A48").Value = 39.79439
$ws.Range("B48").Value = -121.8975
$ws.Range("A49").Value = 27.483056
$ws.Range("B49").Value = -81.920556
$ws.Range("A50").Value = 41.869435
$ws.Range("B50").Value = -94.67733200000001
$ws.Range("A51").Value = 40.904486
$ws.Range("B51").Value = -74.101294
$ws.Range("A52").Value = 45.75538
$ws.Range("B52").Value = -108.531811
$ws.Range("A53").Value = 42.5425
$ws.Range("B53").Value = -88.1704
$ws.Range("A54").Value = 37.32281
$ws.Range("B54").Value = -121.98077
$ws.Range("A55").Value = 42.91601
$ws.Range("B55").Value = -76.79192
$ws.Range("A56").Value = 45.8901
$ws.Range("B56").Value = -88.2796
$ws.Range("A57").Value = 37.35178
$ws.Range("B57").Value = -122.0545
$ws.Range("A58").Value = 44.617384
$ws.Range("B58").Value = -123.103934
$ws.Range("A59").Value = 38.59559
$ws.Range("B59").Value = -75.17868
$ws.Range("A60").Value = 34.02587
$ws.Range("B60").Value = -118.39185
$ws.Range("A61").Value = 45.95746
$ws.Range("B61").Value = -94.66886
$ws.Range("A62").Value = 40.66797
$ws.Range("B62").Value = -92.25228
$ws.Range("A63").Value = 38.58054
$ws.Range("B63").Value = -121.52866
$ws.Range("A64").Value = 38.20591
$ws.Range("B64").Value = -85.66679000000001
$ws.Range("A65").Value = 31.713456
$ws.Range("B65").Value = -82.390362
$ws.Range("A66").Value = 32.32822
$ws.Range("B66").Value = -110.98418
$ws.Range("A67").Value = 37.132856
$ws.Range("B67").Value = -75.966858
$ws.Range("A68").Value = 41.02704
$ws.Range("B68").Value = -73.62445
$ws.Range("A69").Value = 44.196057
$ws.Range("B69").Value = -120.817688
$ws.Range("A70").Value = 40.42475
$ws.Range("B70").Value = -86.88245000000001
$ws.Range("A71").Value = 34.97721
$ws.Range("B71").Value = -109.823346
$ws.Range("A72").Value = 30.09444
$ws.Range("B72").Value = -98.04929
$ws.Range("A73").Value = 45.6079
$ws.Range("B73").Value = -87.94889999999999
$ws.Range("A74").Value = 46.349523
$ws.Range("B74").Value = -85.52015400000001
$ws.Range("A75").Value = 33.59984
$ws.Range("B75").Value = -111.98398
$ws.Range("A76").Value = 47.50617
$ws.Range("B76").Value = -111.25288
$ws.Range("A77").Value = 40.68711
$ws.Range("B77").Value = -73.852755
$ws.Range("A78").Value = 36.788374
$ws.Range("B78").Value = -83.69761
$ws.Range("A79").Value = 40.94577
$ws.Range("B79").Value = -74.28653
$ws.Range("A80").Value = 39.67377
$ws.Range("B80").Value = -74.22311000000001
$ws.Range("A81").Value = 40.85458
$ws.Range("B81").Value = -74.74531
$ws.Range("A82").Value = 40.91078
$ws.Range("B82").Value = -74.40692
$ws.Range("A83").Value = 39.62284
$ws.Range("B83").Value = -74.64192
$ws.Range("A84").Value = 42.26209
$ws.Range("B84").Value = -71.10914
$ws.Range("A85").Value = 43.182449
$ws.Range("B85").Value = -95.856677
$ws.Range("A86").Value = 45.726157
$ws.Range("B86").Value = -108.648217
$ws.Range("A87").Value = 44.1533
$ws.Range("B87").Value = -94.03767999999999
$ws.Range("A88").Value = 35.1754
$ws.Range("B88").Value = -90.0445
$ws.Range("A89").Value = 34.75434
$ws.Range("B89").Value = -84.94982
$ws.Range("A90").Value = 38.496137
$ws.Range("B90").Value = -121.658605
$ws.Range("A91").Value = 32.56683
$ws.Range("B91").Value = -97.42749000000001
$ws.Range("A92").Value = 40.535707
$ws.Range("B92").Value = -74.40082200000001
$ws.Range("A93").Value = 41.55596
$ws.Range("B93").Value = -112.12444
$ws.Range("A94").Value = 30.3294
$ws.Range("B94").Value = -95.46576
$ws.Range("A95").Value = 42.124193
$ws.Range("B95").Value = -72.567819
$ws.Range("A96").Value = 41.94011
$ws.Range("B96").Value = -87.68792000000001
$ws.Range("A97").Value = 38.86383
$ws.Range("B97").Value = -104.81004
$ws.Range("A98").Value = 40.89795
$ws.Range("B98").Value = -72.31458000000001
$ws.Range("A99").Value = 18.440917
$ws.Range("B99").Value = -67.129047
$ws.Range("A100").Value = 39.84214
$ws.Range("B100").Value = -96.64026
$ws.Range("A101").Value = 41.904339
$ws.Range("B101").Value = -88.339822
